# Generate Report for Handoff
# Refresh the handoff report: the four "Ready for handoff" rows (rows 4-7)
# in both the zh-cn and de-de sheets get a new Priority ("ht" instead of
# "low") and an updated "Latest Handoff Datetime" timestamp reflecting the
# new handoff generation run.

$wb = $excel.ActiveWorkbook

$ws_zh = $wb.Worksheets.Item("zh-cn")
$ws_zh.Range("E4:E7").Value = "ht"
$ws_zh.Range("H4:H7").Value = "2016-08-28 12:31:25"

$ws_de = $wb.Worksheets.Item("de-de")
$ws_de.Range("E4:E7").Value = "ht"
$ws_de.Range("H4:H7").Value = "2016-08-28 12:31:31"

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_overview.Range("G4:G7").Value = "2016-08-28 12:31:31"
